$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rnaSampleNumber (column F) for rows 28-53: add 26 to continue numbering
for ($r = 28; $r -le 53; $r++) {
    $old = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 6).Value = $old + 26
}

# Update the selection to reflect the newly entered range (F2:F53), active cell F2
$ws.Range("F2:F53").Select()
